$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cx3cl1"
$ws.Range("C2").Value = "Itgb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 6.662215666666666
$ws.Range("H2").Value = 19.986647
$ws.Range("I2").Value = 0.4964100230945138
$ws.Range("J2").Value = 0.4964100230945138
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.970048
$ws.Range("N2").Value = 26.910144
$ws.Range("O2").Value = 0.487108783009476
$ws.Range("P2").Value = 0.4871087830094759
$ws.Range("Q2").Value = 59.76039431635199
$ws.Range("R2").Value = 537.8435488471679
$ws.Range("S2").Value = 0.2418056822232745
$ws.Range("T2").Value = 0.2418056822232744

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Cx3cl1"
$ws.Range("C3").Value = "Itgb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 6.662215666666666
$ws.Range("H3").Value = 19.986647
$ws.Range("I3").Value = 0.4964100230945138
$ws.Range("J3").Value = 0.4964100230945138
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.012070666666666
$ws.Range("N3").Value = 27.036212
$ws.Range("O3").Value = 0.489390778604016
$ws.Range("P3").Value = 0.489390778604016
$ws.Range("Q3").Value = 60.04035838457376
$ws.Range("R3").Value = 540.3632254611639
$ws.Range("S3").Value = 0.2429384877090617
$ws.Range("T3").Value = 0.2429384877090617

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Cx3cl1"
$ws.Range("C4").Value = "Itgb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 6.662215666666666
$ws.Range("H4").Value = 19.986647
$ws.Range("I4").Value = 0.4964100230945138
$ws.Range("J4").Value = 0.4964100230945138
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.4327576666666667
$ws.Range("N4").Value = 1.298273
$ws.Range("O4").Value = 0.02350043838650813
$ws.Range("P4").Value = 0.02350043838650813
$ws.Range("Q4").Value = 2.883124906736777
$ws.Range("R4").Value = 25.948124160631
$ws.Range("S4").Value = 0.0116658531621777
$ws.Range("T4").Value = 0.0116658531621777

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Cx3cl1"
$ws.Range("C5").Value = "Itgb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 5.098534333333333
$ws.Range("H5").Value = 15.295603
$ws.Range("I5").Value = 0.3798981709375522
$ws.Range("J5").Value = 0.3798981709375522
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.970048
$ws.Range("N5").Value = 26.910144
$ws.Range("O5").Value = 0.487108783009476
$ws.Range("P5").Value = 0.4871087830094759
$ws.Range("Q5").Value = 45.734097699648
$ws.Range("R5").Value = 411.606879296832
$ws.Range("S5").Value = 0.1850517357129169
$ws.Range("T5").Value = 0.1850517357129169

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Cx3cl1"
$ws.Range("C6").Value = "Itgb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 5.098534333333333
$ws.Range("H6").Value = 15.295603
$ws.Range("I6").Value = 0.3798981709375522
$ws.Range("J6").Value = 0.3798981709375522
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.012070666666666
$ws.Range("N6").Value = 27.036212
$ws.Range("O6").Value = 0.489390778604016
$ws.Range("P6").Value = 0.489390778604016
$ws.Range("Q6").Value = 45.94835170842622
$ws.Range("R6").Value = 413.535165375836
$ws.Range("S6").Value = 0.1859186616653702
$ws.Range("T6").Value = 0.1859186616653702

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Cx3cl1"
$ws.Range("C7").Value = "Itgb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 5.098534333333333
$ws.Range("H7").Value = 15.295603
$ws.Range("I7").Value = 0.3798981709375522
$ws.Range("J7").Value = 0.3798981709375522
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.4327576666666667
$ws.Range("N7").Value = 1.298273
$ws.Range("O7").Value = 0.02350043838650813
$ws.Range("P7").Value = 0.02350043838650813
$ws.Range("Q7").Value = 2.206429821513222
$ws.Range("R7").Value = 19.857868393619
$ws.Range("S7").Value = 0.008927773559265079
$ws.Range("T7").Value = 0.008927773559265079

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Cx3cl1"
$ws.Range("C8").Value = "Itgb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.660042
$ws.Range("H8").Value = 4.980126
$ws.Range("I8").Value = 0.1236918059679339
$ws.Range("J8").Value = 0.123691805967934
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 8.970048
$ws.Range("N8").Value = 26.910144
$ws.Range("O8").Value = 0.487108783009476
$ws.Range("P8").Value = 0.4871087830094759
$ws.Range("Q8").Value = 14.890656422016
$ws.Range("R8").Value = 134.015907798144
$ws.Range("S8").Value = 0.06025136507328454
$ws.Range("T8").Value = 0.06025136507328454

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Cx3cl1"
$ws.Range("C9").Value = "Itgb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.660042
$ws.Range("H9").Value = 4.980126
$ws.Range("I9").Value = 0.1236918059679339
$ws.Range("J9").Value = 0.123691805967934
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 9.012070666666666
$ws.Range("N9").Value = 27.036212
$ws.Range("O9").Value = 0.489390778604016
$ws.Range("P9").Value = 0.489390778604016
$ws.Range("Q9").Value = 14.96041581363467
$ws.Range("R9").Value = 134.643742322712
$ws.Range("S9").Value = 0.06053362922958407
$ws.Range("T9").Value = 0.06053362922958407

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Cx3cl1"
$ws.Range("C10").Value = "Itgb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.660042
$ws.Range("H10").Value = 4.980126
$ws.Range("I10").Value = 0.1236918059679339
$ws.Range("J10").Value = 0.123691805967934
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.4327576666666667
$ws.Range("N10").Value = 1.298273
$ws.Range("O10").Value = 0.02350043838650813
$ws.Range("P10").Value = 0.02350043838650813
$ws.Range("Q10").Value = 0.7183959024886667
$ws.Range("R10").Value = 6.465563122398001
$ws.Range("S10").Value = 0.00290681166506535
$ws.Range("T10").Value = 0.00290681166506535
